# Update the "cartao ponto" worksheet: extend the month-marker column (A)
# down through the rows that were previously left blank (101-138), mirroring
# the same pattern already used for the preceding months.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 101-119 continue the "10/2011" (October) marker
$ws.Range("A101:A119").Value = "10/2011"

# Rows 120-138 continue the "11/2011" (November) marker
$ws.Range("A120:A138").Value = "11/2011"
